# Daily attendance processing - 2026-01-05 05:00:07
#
# For the rows below, the "Recorded By" value in column G lists two
# contributors separated by ", ". This pass normalizes the ordering of
# those two names/emails (swapping the first and second entries) for the
# rows identified by the daily processing job.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Rows whose column-G "Recorded By" entry needs its two comma-separated
# parts swapped.
$rowsToSwap = @(
    4, 7, 10, 11, 12, 13, 14, 15, 17, 18, 19, 20, 21, 22, 24, 26,
    30, 33, 36, 37, 38, 39, 40, 41, 43, 44, 45, 46, 47, 48, 50, 52,
    56, 59, 62, 63, 64, 65, 66, 67, 69, 70, 71, 72, 73, 74, 76, 78,
    83, 84, 85, 86, 87, 90, 92, 93, 94, 96, 99, 101,
    109, 110, 111, 112, 113, 116, 118, 119, 120, 122, 125, 127,
    135, 136, 137, 138, 139, 142, 144, 145, 146, 148, 151, 153
)

foreach ($row in $rowsToSwap) {
    $cell = $ws.Range("G$row")
    # NOTE: reading via ".Value" on this host does not return the live
    # scalar (it echoes the property descriptor), so the current text is
    # read through ".Value2" instead; writing back through ".Value" works
    # as expected.
    $current = [string]$cell.Value2
    $parts = $current -split ', ', 2
    if ($parts.Count -eq 2) {
        $cell.Value = "{0}, {1}" -f $parts[1], $parts[0]
    }
}
